$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Friday (F5) and Saturday (G5) hours for the week of 43143.
# I5 holds a shared formula SUM(B5:H5) that will recalc automatically.
$ws.Range("F5").Value = 8.25
$ws.Range("G5").Value = 7.75

# Update the active cell selection to H5 to match the saved view state.
$ws.Range("H5").Select()
